$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I = DAMSLTag (col 9), Column J = DialogAct (col 10)
$updates = @(
    @{ Row = 10; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 20; I = "ba"; J = "Appreciation" }
    @{ Row = 21; I = "sv"; J = "Statement-opinion" }
    @{ Row = 25; I = "ba"; J = "Appreciation" }
    @{ Row = 27; I = "ba"; J = "Appreciation" }
    @{ Row = 28; I = "aa"; J = "Agree/Accept" }
    @{ Row = 33; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 41; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 42; I = "ba"; J = "Appreciation" }
    @{ Row = 43; I = "sv"; J = "Statement-opinion" }
    @{ Row = 45; I = "aa"; J = "Agree/Accept" }
    @{ Row = 56; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 74; I = "aa"; J = "Agree/Accept" }
    @{ Row = 79; I = "ba"; J = "Appreciation" }
    @{ Row = 80; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 87; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 90; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 94; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 106; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 129; I = "%"; J = "Uninterpretable" }
    @{ Row = 135; I = "sv"; J = "Statement-opinion" }
    @{ Row = 141; I = "aa"; J = "Agree/Accept" }
    @{ Row = 152; I = "aa"; J = "Agree/Accept" }
    @{ Row = 157; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 159; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 173; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 174; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 178; I = "ba"; J = "Appreciation" }
    @{ Row = 180; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 184; I = "aa"; J = "Agree/Accept" }
    @{ Row = 186; I = "ba"; J = "Appreciation" }
    @{ Row = 191; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 206; I = "sv"; J = "Statement-opinion" }
    @{ Row = 211; I = "aa"; J = "Agree/Accept" }
    @{ Row = 218; I = "ba"; J = "Appreciation" }
    @{ Row = 229; I = "sv"; J = "Statement-opinion" }
    @{ Row = 236; I = "sv"; J = "Statement-opinion" }
    @{ Row = 256; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 274; I = "sv"; J = "Statement-opinion" }
    @{ Row = 276; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 288; I = "aa"; J = "Agree/Accept" }
    @{ Row = 303; I = "sv"; J = "Statement-opinion" }
    @{ Row = 339; I = "aa"; J = "Agree/Accept" }
    @{ Row = 340; I = "aa"; J = "Agree/Accept" }
    @{ Row = 349; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 352; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 355; I = "aa"; J = "Agree/Accept" }
    @{ Row = 418; I = "%"; J = "Uninterpretable" }
    @{ Row = 425; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 428; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 432; I = "sv"; J = "Statement-opinion" }
    @{ Row = 440; I = "ba"; J = "Appreciation" }
    @{ Row = 441; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 459; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 461; I = "aa"; J = "Agree/Accept" }
    @{ Row = 462; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 464; I = "b"; J = "Acknowledge (Backchannel)" }
    @{ Row = 472; I = "%"; J = "Uninterpretable" }
    @{ Row = 478; I = "sd"; J = "Statement-non-opinion" }
    @{ Row = 483; I = "aa"; J = "Agree/Accept" }
    @{ Row = 512; I = "b"; J = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"